$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Uur=2, project versie=1.1, Beschrijving=<long text> ---

# A3: plain number
$ws.Range("A3").Value = 2

# B3: copy formatting (style) from B2, then write "1.1" as literal text
# (not a number). A direct .Value = "1.1" assignment gets auto-coerced to
# the number 1.1 by the engine, which would store it as a <v> numeric
# cell instead of a shared-string text cell. Routing the text through a
# formula that evaluates to a string, then pasting values-only into B3,
# keeps it a literal string without creating any extra number-format
# style in styles.xml.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("ZZ1").Formula = '="1.1"'
$ws.Range("ZZ1").Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("ZZ1").Clear()

# C3: copy formatting from C2, then set the long description text
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C3").Value = 'In dit uur heb ik een in en uit zoom functie script gemaakt voor de camera. Zo kan de user makkelijker de hele maze zien. Daarna heb ik de code geschreven om de hoogte van de muren te kunen instellen. Als laatste heb ik ervoor gezorgd dat de user een nieuwe maze kan generaten door op "G" te drukken. De oude maze word hierdoor verwijderd en daarna word de nieuwe gemaakt.'

$excel.CutCopyMode = 0

# --- Column C width (194.42578125 -> 255.7109375 chars) ---
$ws.Columns("C").ColumnWidth = 254.8333333333

# --- View: scroll so column C is the leftmost visible column, select C6 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C6").Select()
